$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Insert a new column before column K (so it becomes the new "County" column
# right after "lon", pushing "Created By (Lookup)" and everything following
# it one column to the right).
$ws.Range("K1").EntireColumn.Insert()

# New column inherits the width used by the neighboring lat/lon columns.
$ws.Range("K1").ColumnWidth = $ws.Range("J1").ColumnWidth

# Header for the new column.
$ws.Range("K1").Value = "County"

# County values for each office location row.
# Order of entry matters for the shared-strings table: new unique strings
# get appended in the order they are first written, so match the original
# authoring order (Fulton, Travis, Baton Rouge, Cool, Boulder, Suffolk).
$ws.Range("K2").Value = "Fulton"
$ws.Range("K3").Value = "Travis"
$ws.Range("K4").Value = "Baton Rouge"
$ws.Range("K7").Value = "Cool"
$ws.Range("K6").Value = "Boulder"
$ws.Range("K5").Value = "Suffolk"

# Grow the table definition to include the new column, then re-touch every
# header cell after it so the table's column list re-syncs its names to the
# (now shifted) worksheet headers.
$tbl.Resize($ws.Range("A1:AG29"))
for ($col = 12; $col -le 33; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $cell.Text
}

$ws.Range("K6").Select()
